# Auto-generated edit script: updates H/I/J/K/L/M/N market-price/profit
# columns on specific Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ================= Sheet ALC =================
$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 1378706.1
$ws.Range("J17").Value = 1421750.8
$ws.Range("L17").Value = 4265252.4
$ws.Range("N17").Value = -4265588.4
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 2997.4
$ws.Range("I62").Value = 3220.4443
$ws.Range("K62").Value = 3220.4443
$ws.Range("M62").Value = -2596.4443
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 2997.4
$ws.Range("I65").Value = 3220.4443
$ws.Range("K65").Value = 16102.2215
$ws.Range("M65").Value = -12982.2215
# Row 94: Magic Beans
$ws.Range("H94").Value = 55689136
$ws.Range("J94").Value = 600000
$ws.Range("L94").Value = 600000
$ws.Range("N94").Value = -600902
# Row 113: Amaro Kart
$ws.Range("H113").Value = 21374.75
$ws.Range("I113").Value = 21249.5
$ws.Range("K113").Value = 21249.5
$ws.Range("M113").Value = -17995.5

# ================= Sheet ARM =================
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 5844.315
$ws.Range("I32").Value = 5736.7393
$ws.Range("J32").Value = 7700
$ws.Range("K32").Value = 5736.7393
$ws.Range("L32").Value = 7700
$ws.Range("M32").Value = -5449.7393
$ws.Range("N32").Value = -8274
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 16248.625
$ws.Range("I74").Value = 22497.6
$ws.Range("K74").Value = 22497.6
$ws.Range("M74").Value = -21623.6
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 16248.625
$ws.Range("I77").Value = 22497.6
$ws.Range("K77").Value = 112488
$ws.Range("M77").Value = -108120
# Row 97: Ore for Me
$ws.Range("H97").Value = 4881013
$ws.Range("I97").Value = 3711.0322
$ws.Range("J97").Value = 20000650
$ws.Range("K97").Value = 3711.0322
$ws.Range("L97").Value = 20000650
$ws.Range("M97").Value = -3215.0322
$ws.Range("N97").Value = -20001642
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4112.24
$ws.Range("I132").Value = 4014.85
$ws.Range("J132").Value = 4501.8
$ws.Range("K132").Value = 12044.55
$ws.Range("L132").Value = 13505.4
$ws.Range("M132").Value = -9514.549999999999
$ws.Range("N132").Value = -18565.4

# ================= Sheet BSM =================
$ws = $wb.Worksheets.Item("BSM")
# Row 38: The Naked Blade
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("N38").ClearContents()
# Row 46: Spice Cadet
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("N46").ClearContents()
# Row 74: I Could Feel That from Here
$ws.Range("H74").Value = 19698.8
$ws.Range("J74").Value = 19698.8
$ws.Range("L74").Value = 19698.8
$ws.Range("N74").Value = -21570.8
# Row 77: Tensions in Creasing (L)
$ws.Range("H77").Value = 19698.8
$ws.Range("J77").Value = 19698.8
$ws.Range("L77").Value = 59096.39999999999
$ws.Range("N77").Value = -68456.39999999999
# Row 87: Winter Weather Conditions
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90: The Nightsoil Is Dark and Full of Terrors (L)
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").ClearContents()
# Row 94: High Steal
$ws.Range("H94").Value = 9341.088
$ws.Range("I94").Value = 11937.625
$ws.Range("K94").Value = 11937.625
$ws.Range("M94").Value = -11486.625
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 10207.156
$ws.Range("I99").Value = 10662.423
$ws.Range("K99").Value = 10662.423
$ws.Range("M99").Value = -9164.423000000001
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 9602.5
$ws.Range("I134").Value = 10042.733
$ws.Range("K134").Value = 30128.199
$ws.Range("M134").Value = -27593.199
# Row 140: Ceremonial Teeth
$ws.Range("H140").Value = 72370.5
$ws.Range("J140").Value = 72370.5
$ws.Range("L140").Value = 72370.5
$ws.Range("N140").Value = -82730.5

# ================= Sheet CRP =================
$ws = $wb.Worksheets.Item("CRP")
# Row 35: Storm of Swords
$ws.Range("H35").Value = 4600
$ws.Range("J35").Value = 4200
$ws.Range("L35").Value = 4200
$ws.Range("N35").Value = -4788
# Row 39: An Expected Tourney
$ws.Range("H39").Value = 9500
$ws.Range("I39").Value = 9500
$ws.Range("K39").Value = 9500
$ws.Range("M39").Value = -9109
# Row 47: Grippy When Wet
$ws.Range("H47").Value = 33717.5
$ws.Range("I47").Value = 20500
$ws.Range("K47").Value = 20500
$ws.Range("M47").Value = -19934
# Row 49: Bend It Like Durendaire
$ws.Range("H49").Value = 9500
$ws.Range("I49").Value = 9500
$ws.Range("K49").Value = 9500
$ws.Range("M49").Value = -9318
# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 13343.6
$ws.Range("I122").Value = 13343.6
$ws.Range("K122").Value = 40030.8
$ws.Range("M122").Value = -37580.8
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 6909.45
$ws.Range("I134").Value = 7961.875
$ws.Range("K134").Value = 23885.625
$ws.Range("M134").Value = -21350.625

# ================= Sheet CUL =================
$ws = $wb.Worksheets.Item("CUL")
# Row 96: Hunger Is No Game
$ws.Range("H96").Value = 4000
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 12000
$ws.Range("N96").Value = -16118
# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 1861
$ws.Range("I98").Value = 1874
$ws.Range("K98").Value = 5622
$ws.Range("M98").Value = -4124
# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 10991762
$ws.Range("J131").Value = 1966.2759
$ws.Range("L131").Value = 5898.8277
$ws.Range("N131").Value = -15978.8277

# ================= Sheet GSM =================
$ws = $wb.Worksheets.Item("GSM")
# Row 46: Burning the Midnight Oil
$ws.Range("H46").Value = 48998.668
$ws.Range("J46").Value = 48998.668
$ws.Range("L46").Value = 48998.668
$ws.Range("N46").Value = -49310.668
# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 5205.5
$ws.Range("I97").Value = 6297.8
$ws.Range("J97").Value = 2474.75
$ws.Range("K97").Value = 6297.8
$ws.Range("L97").Value = 2474.75
$ws.Range("M97").Value = -5801.8
$ws.Range("N97").Value = -3466.75
# Row 126: Gold Rush Order
$ws.Range("H126").Value = 8406.637000000001
$ws.Range("I126").Value = 11291.308
$ws.Range("K126").Value = 33873.924
$ws.Range("M126").Value = -31403.924

# ================= Sheet LTW =================
$ws = $wb.Worksheets.Item("LTW")
# Row 8: Mind over Muzzle
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("N8").ClearContents()
# Row 16: Saddle Sore
$ws.Range("H16").Value = 3652.3125
$ws.Range("I16").Value = 3906.6924
$ws.Range("J16").Value = 2550
$ws.Range("K16").Value = 3906.6924
$ws.Range("L16").Value = 2550
$ws.Range("M16").Value = -3736.6924
$ws.Range("N16").Value = -2890
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 7745585
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 13941413
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 13941413
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -13941789
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 4586.4
$ws.Range("I61").Value = 1882.7333
$ws.Range("J61").Value = 12697.4
$ws.Range("K61").Value = 1882.7333
$ws.Range("L61").Value = 12697.4
$ws.Range("M61").Value = -1680.7333
$ws.Range("N61").Value = -13101.4
# Row 62: Pummeling Abroad
$ws.Range("H62").Value = 15000
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16248
# Row 65: The Style of the Time (L)
$ws.Range("H65").Value = 15000
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51240
# Row 113: Peace in Rest
$ws.Range("H113").Value = 4586.4
$ws.Range("I113").Value = 1882.7333
$ws.Range("J113").Value = 12697.4
$ws.Range("K113").Value = 1882.7333
$ws.Range("L113").Value = 12697.4
$ws.Range("M113").Value = 287.2666999999999
$ws.Range("N113").Value = -17037.4
# Row 122: Hell on Leather
$ws.Range("H122").Value = 5350.9033
$ws.Range("I122").Value = 5250.6665
$ws.Range("J122").Value = 5561.4
$ws.Range("K122").Value = 15751.9995
$ws.Range("L122").Value = 16684.2
$ws.Range("M122").Value = -13301.9995
$ws.Range("N122").Value = -21584.2
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 995656.1
$ws.Range("I132").Value = 1355986.1
$ws.Range("J132").Value = 4748.75
$ws.Range("K132").Value = 4067958.3
$ws.Range("L132").Value = 14246.25
$ws.Range("M132").Value = -4065428.3
$ws.Range("N132").Value = -19306.25

# ================= Sheet WVR =================
$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 289346.9
$ws.Range("I62").Value = 569164.7
$ws.Range("J62").Value = 9529.166999999999
$ws.Range("K62").Value = 569164.7
$ws.Range("L62").Value = 9529.166999999999
$ws.Range("M62").Value = -568540.7
$ws.Range("N62").Value = -10777.167
# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 289346.9
$ws.Range("I65").Value = 569164.7
$ws.Range("J65").Value = 9529.166999999999
$ws.Range("K65").Value = 2845823.5
$ws.Range("L65").Value = 47645.835
$ws.Range("M65").Value = -2842703.5
$ws.Range("N65").Value = -53885.835
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 16747.555
$ws.Range("I126").Value = 21619.3
$ws.Range("J126").Value = 2828.2856
$ws.Range("K126").Value = 64857.89999999999
$ws.Range("L126").Value = 8484.856800000001
$ws.Range("M126").Value = -62387.89999999999
$ws.Range("N126").Value = -13424.8568
